# Applies the betexplorer San Marino "campionato-sammarinese" 2023-2024 update
# (script run 12-11-2023 14:45): corrects the home/away ordering for several
# already-recorded fixtures and appends three newly played fixtures (rows 63-65).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Fix rows whose match data got shuffled between sibling rows ---------
$rowData = @{
  5 = @{ F='Cailungo'; G=0; H='Folgore'; I=2; J=4.61; L=3.25; M='16/09/2023 14:24'; N=3.88; P=3.77; Q='16/09/2023 14:24'; R=1.48; T=1.88; U='16/09/2023 14:24'; V='https://www.betexplorer.com/football/san-marino/campionato-sammarinese/cailungo-folgore/8W6ZyDI3/' }
  6 = @{ F='Murata'; G=5; H='Faetano'; I=0; J=1.95; L=1.54; M='16/09/2023 14:56'; N=3.35; P=4.83; Q='16/09/2023 14:56'; R=2.98; T=4.02; U='16/09/2023 14:56'; V='https://www.betexplorer.com/football/san-marino/campionato-sammarinese/ss-murata-sc-faetano/C4DQwZmi/' }
  7 = @{ F='Tre Fiori'; G=0; H='Virtus'; I=1; J=2.39; L=2.49; M='16/09/2023 12:06'; N=2.85; P=3.09; Q='16/09/2023 13:02'; R=2.65; T=2.6; U='16/09/2023 12:06'; V='https://www.betexplorer.com/football/san-marino/campionato-sammarinese/tre-fiori-virtus/b1PrFkBj/' }
  16 = @{ F='Tre Fiori'; G=1; H='Fiorentino'; I=0; J=1.44; L=1.5; M='30/09/2023 14:01'; N=3.82; P=4.21; Q='30/09/2023 14:06'; R=5.1; T=4.99; U='30/09/2023 14:06'; V='https://www.betexplorer.com/football/san-marino/campionato-sammarinese/tre-fiori-fiorentino/byQ00oPc/' }
  19 = @{ F='Cosmos'; G=0; H='San Giovanni'; I=0; J=1.15; L=1.13; M='30/09/2023 11:15'; N=5.91; P=6.92; Q='30/09/2023 13:02'; R=9.56; T=12.95; U='30/09/2023 11:15'; V='https://www.betexplorer.com/football/san-marino/campionato-sammarinese/sp-cosmos-san-giovanni/6swFxrXj/' }
  31 = @{ F='Tre Fiori'; G=3; H='Folgore'; I=1; J=1.63; L=1.63; M='07/10/2023 12:14'; N=3.34; P=3.65; Q='07/10/2023 13:02'; R=4.24; T=4.51; U='07/10/2023 12:14'; V='https://www.betexplorer.com/football/san-marino/campionato-sammarinese/tre-fiori-folgore/84SBAj42/' }
  32 = @{ F='Tre Penne'; G=5; H='San Giovanni'; I=0; J=1.12; L=1.12; M='07/10/2023 14:12'; N=6.37; P=7.58; Q='07/10/2023 14:12'; R=10.28; T=12.66; U='07/10/2023 14:12'; V='https://www.betexplorer.com/football/san-marino/campionato-sammarinese/tre-penne-san-giovanni/6unHl9sq/' }
  45 = @{ F='San Giovanni'; G=2; H='Fiorentino'; I=0; J=3.06; L=3.24; M='28/10/2023 14:33'; N=3.21; P=3.24; Q='28/10/2023 14:33'; R=1.96; T=2.05; U='28/10/2023 14:33'; V='https://www.betexplorer.com/football/san-marino/campionato-sammarinese/san-giovanni-fiorentino/xba1L7IK/' }
  46 = @{ F='Pennarossa'; G=1; H='Domagnano'; I=2; J=2.35; L=3.25; M='28/10/2023 14:46'; N=3.1; P=3.66; Q='28/10/2023 14:33'; R=2.56; T=1.9; U='28/10/2023 14:46'; V='https://www.betexplorer.com/football/san-marino/campionato-sammarinese/ss-pennarossa-sp-domagnano/UevtEPAf/' }
  48 = @{ F='Tre Fiori'; G=5; H='Cosmos'; I=0; J=3.11; L=2.84; M='29/10/2023 09:05'; N=2.97; P=3.25; Q='29/10/2023 13:01'; R=2.05; T=2.2; U='29/10/2023 09:05'; V='https://www.betexplorer.com/football/san-marino/campionato-sammarinese/tre-fiori-sp-cosmos/8OfcMm3E/' }
  49 = @{ F='Cailungo'; G=4; H='Faetano'; I=1; J=2.69; L=3.01; M='29/10/2023 14:54'; N=3.3; P=3.7; Q='29/10/2023 14:54'; R=2.12; T=1.99; U='29/10/2023 14:54'; V='https://www.betexplorer.com/football/san-marino/campionato-sammarinese/cailungo-sc-faetano/MVzYFotr/' }
  50 = @{ F='Tre Penne'; G=3; H='Murata'; I=0; J=1.45; L=1.85; M='29/10/2023 14:24'; N=4.07; P=3.81; Q='29/10/2023 14:32'; R=4.66; T=3.26; U='29/10/2023 14:24'; V='https://www.betexplorer.com/football/san-marino/campionato-sammarinese/tre-penne-ss-murata/zZvxF5el/' }
  56 = @{ F='Tre Fiori'; G=5; H='Juvenes/Dogana'; I=1; J=1.44; L=1.56; M='05/11/2023 14:08'; N=3.9; P=4.05; Q='05/11/2023 14:08'; R=5.01; T=4.62; U='05/11/2023 14:08'; V='https://www.betexplorer.com/football/san-marino/campionato-sammarinese/tre-fiori-juvenes-dogana/M91DEd15/' }
  57 = @{ F='Murata'; G=5; H='Pennarossa'; I=0; J=1.37; L=1.39; M='05/11/2023 14:06'; N=4.38; P=4.73; Q='05/11/2023 14:06'; R=5.5; T=5.7; U='05/11/2023 14:06'; V='https://www.betexplorer.com/football/san-marino/campionato-sammarinese/ss-murata-ss-pennarossa/4YjIDGGB/' }
}

foreach ($r in $rowData.Keys) {
    $data = $rowData[$r]
    foreach ($col in $data.Keys) {
        $ws.Range("$col$r").Value = $data[$col]
    }
}

# --- 2) Append the three newly scraped fixtures (rows 63-65) ---------------
$newRowData = @{
  63 = @{ A=62; B='san-marino'; C='campionato-sammarinese'; D='2023-2024'; E=45242.625; F='Libertas'; G=1; H='San Giovanni'; I=1; J=1.5; K='12/11/2023 04:12'; L=1.47; M='12/11/2023 14:06'; N=3.96; O='12/11/2023 04:12'; P=4.35; Q='12/11/2023 14:48'; R=4.89; S='12/11/2023 04:12'; T=5.13; U='12/11/2023 14:06'; V='https://www.betexplorer.com/football/san-marino/campionato-sammarinese/ac-libertas-san-giovanni/IyHiSFV4/' }
  64 = @{ A=63; B='san-marino'; C='campionato-sammarinese'; D='2023-2024'; E=45242.625; F='Juvenes/Dogana'; G=0; H='Cosmos'; I=4; J=5.4; K='12/11/2023 04:12'; L=7.74; M='12/11/2023 13:32'; N=4.28; O='12/11/2023 04:12'; P=4.97; Q='12/11/2023 13:32'; R=1.42; S='12/11/2023 04:12'; T=1.29; U='12/11/2023 13:32'; V='https://www.betexplorer.com/football/san-marino/campionato-sammarinese/juvenes-dogana-sp-cosmos/M9NrUy1h/' }
  65 = @{ A=64; B='san-marino'; C='campionato-sammarinese'; D='2023-2024'; E=45242.625; F='Virtus'; G=6; H='Faetano'; I=1; J=1.21; K='12/11/2023 04:12'; L=1.19; M='12/11/2023 14:53'; N=5.7; O='12/11/2023 04:12'; P=6.22; Q='12/11/2023 14:53'; R=8.26; S='12/11/2023 04:12'; T=9.65; U='12/11/2023 14:53'; V='https://www.betexplorer.com/football/san-marino/campionato-sammarinese/virtus-sc-faetano/lMhUAE0U/' }
}

$lastDataRow = 62

foreach ($r in ($newRowData.Keys | Sort-Object)) {
    # Copy the formatting (bold index column, date format, borders) from the
    # last existing data row so the new rows look like the rest of the table.
    $ws.Range("A" + $lastDataRow + ":V" + $lastDataRow).Copy()
    $ws.Range("A" + $r + ":V" + $r).PasteSpecial(-4122)

    $data = $newRowData[$r]
    foreach ($col in $data.Keys) {
        $ws.Range("$col$r").Value = $data[$col]
    }
}

$excel.CutCopyMode = 0
